# header, footer, rendered, exportable周り実装
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the TODO placeholders to their real template tokens.
$ws.Range("A1").Value = "`$C[]{headers}"
$ws.Range("A3").Value = "`$C[]{footers}"

# A1 (header) and A3 (footer) get a thin box border plus a light
# accent-5 fill; A2 (data columns) gets the same border but no fill.
$headerFooter = $ws.Range("A1:A1")
$headerFooter2 = $ws.Range("A3:A3")
$dataRow = $ws.Range("A2:A2")

foreach ($rng in @($headerFooter, $dataRow, $headerFooter2)) {
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
}

$headerFooter.Interior.ThemeColor = 8
$headerFooter.Interior.TintAndShade = 0.79998168889431442

$headerFooter2.Interior.ThemeColor = 8
$headerFooter2.Interior.TintAndShade = 0.79998168889431442
